$wb = $excel.ActiveWorkbook

# Rows to update: row number -> new F value (or $null if unchanged)
$updates = @{
    2  = 159
    3  = 1733
    4  = $null
    5  = $null
    6  = $null
    7  = 12040
    8  = $null
    9  = $null
    10 = 482
    11 = 418
    12 = $null
    13 = 872
    14 = 13497
    15 = 13529
    16 = $null
    17 = $null
    18 = $null
    19 = $null
    20 = 990
    21 = $null
    22 = $null
    23 = 1938
    24 = $null
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 24; $row++) {
        $eCell = $ws.Cells.Item($row, 5)   # column E
        $eCell.Value = $eCell.Value2 -replace '-(\d\d\.\d\d)', ' - $1'

        $newF = $updates[$row]
        if ($null -ne $newF) {
            $ws.Cells.Item($row, 6).Value = $newF
        }
    }
}
